# Add Batteries in addition to Ultracaps
# - gStation sheet: insert a new row above "ultracap.p" for "elecSto_type"
#   and two new rows after "ultracap.N" for "batt.p" / "batt.N".
# - gStation becomes the active/selected sheet (was "tether").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gStation")

# --- Insert "elecSto_type" row right above the existing "ultracap.p" row (row 20) ---
$ws.Rows("20:20").Insert()
$ws.Cells.Item(20, 1).Value2 = "elecSto_type"
$ws.Cells.Item(20, 2).Value2 = 1

# "ultracap.p" / "ultracap.N" are now on rows 21 / 22 (shifted down by the insert above).

# --- Insert two rows for the battery entries right after "ultracap.N" (row 22) ---
$ws.Rows("23:24").Insert()

$ws.Cells.Item(23, 1).Value2 = "batt.p"
$ws.Cells.Item(23, 2).Value2 = 180
$ws.Cells.Item(23, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(24, 1).Value2 = "batt.N"
$ws.Cells.Item(24, 2).Formula = "=10000"

# --- Make "gStation" the selected / active sheet, with G21 as the active cell ---
$ws.Activate()
$ws.Range("G21").Select()
